# Add two new columns (I: "I0", J: "IF") to the worksheet, matching the
# header style of the existing header row (copy formats from H1, which
# carries the bold/bordered/centered header style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2-29 in columns I and J.
$values = @{
    2  = @(8, 8)
    3  = @(9, 10)
    4  = @(9, 9)
    5  = @(9, 9)
    6  = @(8, 9)
    7  = @(9, 9)
    8  = @(8, 8)
    9  = @(7, 7)
    10 = @(8, 8)
    11 = @(9, 9)
    12 = @(10, 10)
    13 = @(8, 8)
    14 = @(8, 8)
    15 = @(5, 5)
    16 = @(7, 7)
    17 = @(10, 10)
    18 = @(5, 6)
    19 = @(6, 6)
    20 = @(6, 6)
    21 = @(6, 6)
    22 = @(6, 6)
    23 = @(4, 4)
    24 = @(6, 6)
    25 = @(5, 6)
    26 = @(5, 5)
    27 = @(7, 7)
    28 = @(8, 8)
    29 = @(4, 4)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
